$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

function Get-ParaByText($text) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Trim() -eq $text) {
            return $p
        }
    }
    return $null
}

# ------------------------------------------------------------------
# 1) "Los registros de gasto se presentaran en forma de grilla"
#    Merge the 3 runs (with gramStart/gramEnd proofErr around
#    "presentaran") into a single plain run.
# ------------------------------------------------------------------
$p1 = Get-ParaByText("Los registros de gasto se presentaran en forma de grilla")
$xml1 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="000D44EC" w:rsidRDefault="00634CED" w:rsidP="000D44EC"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Los registros de gasto se presentaran en forma de grilla</w:t></w:r></w:p>'
$p1.Range.InsertXML($xml1)

# ------------------------------------------------------------------
# 2) New bullet after "En la casilla de nombre al ir digitando..."
#    "No debe permitir ingresar movimientos en meses ya cerrados."
# ------------------------------------------------------------------
$p2 = Get-ParaByText("En la casilla de nombre al ir digitando me ayude con sugerencias a partir de gastos ya registrados.")
$p2.Range.InsertParagraphAfter()
$p2new = Get-ParaByText("")
# locate the freshly inserted empty paragraph: it is the paragraph
# right after $p2 in document order
$p2new = $p2.Next()
$xml2 = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>No debe permitir ingresar movimientos en meses ya cerrados.</w:t></w:r></w:p>'
$p2new.Range.InsertXML($xml2)

# ------------------------------------------------------------------
# 3)/4) Move <w:lastRenderedPageBreak/> from the "Sumatoria de los
#        precios" run to the "Fecha" run right before it.
# ------------------------------------------------------------------
$p3 = Get-ParaByText("Fecha")
$xml3 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00C377C8" w:rsidRDefault="00634CED" w:rsidP="00C377C8"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Fecha</w:t></w:r></w:p>'
$p3.Range.InsertXML($xml3)

$p4 = Get-ParaByText("Sumatoria de los precios")
$xml4 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00C377C8" w:rsidRDefault="00C377C8" w:rsidP="00C377C8"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Sumatoria de los precios</w:t></w:r></w:p>'
$p4.Range.InsertXML($xml4)

# ------------------------------------------------------------------
# 5) "... cuanto gasto se presentara en el mes consultado, teniendo
#    en cuenta:" - merge the runs around "presentara" (gramStart/End)
# ------------------------------------------------------------------
$p5 = Get-ParaByText("Tomando la información de meses previos preveer cuanto gasto se presentara en el mes consultado, teniendo en cuenta:")
$xml5 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00C377C8" w:rsidRDefault="00C377C8" w:rsidP="00C377C8"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Tomando la información de meses previos </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>preveer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> cuanto gasto se presentara en el mes consultado, teniendo en cuenta:</w:t></w:r></w:p>'
$p5.Range.InsertXML($xml5)

# ------------------------------------------------------------------
# 6)/7) Move the _GoBack bookmark from the end of "Paramétrica de
#        categorías de gastos:" paragraph to the end of the
#        "Comportamiento de ese mismo mes..." paragraph.
# ------------------------------------------------------------------
$p6 = Get-ParaByText("Comportamiento de ese mismo mes en años anteriores con respecto a otros meses")
$xml6 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00C377C8" w:rsidRDefault="00C377C8" w:rsidP="00C377C8"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t>Comportamiento de ese mismo mes en años anteriores con respecto a otros meses</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$p6.Range.InsertXML($xml6)

$p7 = Get-ParaByText("Paramétrica de categorías de gastos:")
$xml7 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00C377C8" w:rsidRDefault="00C377C8" w:rsidP="00C377C8"><w:r><w:t>Paramétrica de categorías de gastos:</w:t></w:r></w:p>'
$p7.Range.InsertXML($xml7)

# ------------------------------------------------------------------
# 8) New bullet after "Al final de cada mes me envié por correo mi
#    balance del mes."
#    "Cierra el periodo del mes para no permitir mas movimientos en
#    el mes"
# ------------------------------------------------------------------
$p8 = Get-ParaByText("Al final de cada mes me envié por correo mi balance del mes.")
$p8.Range.InsertParagraphAfter()
$p8new = $p8.Next()
$xml8 = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Cierra el periodo del mes para no permitir </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> movimientos en el mes</w:t></w:r></w:p>'
$p8new.Range.InsertXML($xml8)

# ------------------------------------------------------------------
# 9) "pdfmake,wkhtmltopdf" -> split into "pdfmake" + ",wkhtmltopdf"
#    runs, moving the gramStart proofErr marker in between.
# ------------------------------------------------------------------
$p9 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*pdfmake*") {
        $p9 = $p
    }
}
$xml9 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00A37A7C" w:rsidRDefault="00ED1541" w:rsidP="000D44EC"><w:r><w:t xml:space="preserve">Opciones </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pdf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> -&gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pdfmake</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>,wkhtmltopdf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>'
$p9.Range.InsertXML($xml9)

Write-Host "Done."
